$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 59: split "btm"/"bottom" into two distinct entries
$ws.Range("B59").Value = "bottom up"
$ws.Range("A59").Value = "BTM UP"

# Add the new row 60 with the remaining pairing
$ws.Range("A60").Value = "BTM"
$ws.Range("B60").Value = "bottom"

# Expand the table (ListObject) to include the newly added row
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A1:B60"))

# Match the final cursor position left behind by the editing session
[void]$ws.Range("B66").Select()
